$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value is unambiguously numeric-looking need the
# NumberFormat forced to text ("@") first, otherwise Excel COM auto-converts
# the assigned string into a numeric cell instead of keeping it as text.
$textForceCells = @(
    "D5"
    "D6"
    "D7"
    "D9"
    "D13"
    "D14"
    "D19"
    "D22"
    "D23"
    "D24"
    "D25"
    "D27"
    "D30"
    "D32"
    "D34"
    "D35"
    "D36"
    "D37"
    "D38"
    "D39"
    "D42"
    "D43"
    "D45"
    "D46"
    "D47"
    "D48"
)
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated cell values (coin list refreshed from the source feed).
$ws.Range("D2").Value = "65.347.88"
$ws.Range("E2").Value = "  +2.01%  "
$ws.Range("D3").Value = "3.189.54"
$ws.Range("E3").Value = "  +5.27%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "569.22"
$ws.Range("E5").Value = "  +3.98%  "
$ws.Range("D6").Value = "147.95"
$ws.Range("E6").Value = "  +7.84%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "3.184.36"
$ws.Range("E8").Value = "  +5.35%  "
$ws.Range("D9").Value = "0.506"
$ws.Range("E9").Value = "  +4.57%  "
$ws.Range("E10").Value = "  +9.18%  "
$ws.Range("E11").Value = "  +5.41%  "
$ws.Range("E12").Value = "  +5.71%  "
$ws.Range("D13").Value = "37.99"
$ws.Range("E13").Value = "  +7.12%  "
$ws.Range("D14").Value = "0.0000229"
$ws.Range("E14").Value = "  +5.35%  "
$ws.Range("D15").Value = "3.700.38"
$ws.Range("E15").Value = "  +5.48%  "
$ws.Range("D16").Value = "65.465.59"
$ws.Range("E16").Value = "  +2.14%  "
$ws.Range("D17").Value = "3.194.26"
$ws.Range("E17").Value = "  +5.60%  "
$ws.Range("E18").Value = "  +2.68%  "
$ws.Range("D19").Value = "529.24"
$ws.Range("E19").Value = "  +10.67%  "
$ws.Range("E20").Value = "  +7.70%  "
$ws.Range("E21").Value = "  +6.35%  "
$ws.Range("D22").Value = "0.733"
$ws.Range("E22").Value = "  +8.41%  "
$ws.Range("D23").Value = "7.67"
$ws.Range("E23").Value = "  +9.33%  "
$ws.Range("D24").Value = "13.19"
$ws.Range("E24").Value = "  +6.88%  "
$ws.Range("D25").Value = "80.29"
$ws.Range("E25").Value = "  +3.24%  "
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "9.23"
$ws.Range("E27").Value = "  +20.58%  "
$ws.Range("E28").Value = "  +7.93%  "
$ws.Range("E29").Value = "  +7.97%  "
$ws.Range("D30").Value = "27.13"
$ws.Range("E30").Value = "  +6.28%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("D32").Value = "2.67"
$ws.Range("E32").Value = "  +1.55%  "
$ws.Range("E33").Value = "  +5.12%  "
$ws.Range("D34").Value = "553.68"
$ws.Range("E34").Value = "  -4.42%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").Value = "6.25"
$ws.Range("E35").Value = "  +8.23%  "
$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.53"
$ws.Range("E36").Value = "  +4.29%  "
$ws.Range("D37").Value = "54.61"
$ws.Range("E37").Value = "  +5.78%  "
$ws.Range("D38").Value = "0.0442"
$ws.Range("E38").Value = "  +7.81%  "
$ws.Range("D39").Value = "0.0845"
$ws.Range("E39").Value = "  +7.44%  "
$ws.Range("E40").Value = "  +6.99%  "
$ws.Range("D41").Value = "3.196.47"
$ws.Range("E41").Value = "  +10.01%  "
$ws.Range("D42").Value = "2.85"
$ws.Range("E42").Value = "  +3.65%  "
$ws.Range("D43").Value = "8.51"
$ws.Range("E43").Value = "  +4.47%  "
$ws.Range("E44").Value = "  +15.09%  "
$ws.Range("D45").Value = "2.28"
$ws.Range("E45").Value = "  +11.27%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "26.15"
$ws.Range("E47").Value = "  +6.19%  "
$ws.Range("D48").Value = "123.58"
$ws.Range("E48").Value = "  +4.56%  "
$ws.Range("E49").Value = "  +2.65%  "
$ws.Range("E50").Value = "  +3.39%  "
$ws.Range("E51").Value = "  +7.40%  "
